$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Add the new 2020 seasonal observation row (row 42), continuing the
# existing year/mean_chla/mean_colour/mean_tp series.
$ws.Cells.Item(42, 1).Value = 2020
$ws.Cells.Item(42, 2).Value = 31.83
$ws.Cells.Item(42, 3).Value = 11.775
$ws.Cells.Item(42, 4).Value = 49.33
